# Update cryptocurrency price/volume data (cryptos list refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'51.545.01"
$ws.Range("E2").Value = "  -0.61%  "

$ws.Range("D3").Value = "'2.938.02"
$ws.Range("E3").Value = "  +0.55%  "

$ws.Range("E4").Value = "  -0.07%  "

$ws.Range("D5").Value = "'357.86"
$ws.Range("E5").Value = "  +0.00%  "

$ws.Range("D6").Value = "'104.86"
$ws.Range("E6").Value = "  -4.16%  "

$ws.Range("E7").Value = "  -3.42%  "

$ws.Range("E8").Value = "  +0.00%  "

$ws.Range("D9").Value = "'0.595"
$ws.Range("E9").Value = "  -5.54%  "

$ws.Range("D10").Value = "'37.35"
$ws.Range("E10").Value = "  -4.83%  "

$ws.Range("E11").Value = "  +2.28%  "

$ws.Range("D12").Value = "'0.0843"
$ws.Range("E12").Value = "  -3.85%  "

$ws.Range("D13").Value = "'18.73"
$ws.Range("E13").Value = "  -4.43%  "

$ws.Range("D14").Value = "'3.400.81"
$ws.Range("E14").Value = "  +0.55%  "

$ws.Range("D15").Value = "'7.42"
$ws.Range("E15").Value = "  -6.12%  "

$ws.Range("D16").Value = "'2.937.54"
$ws.Range("E16").Value = "  +0.87%  "

$ws.Range("D17").Value = "'0.980"
$ws.Range("E17").Value = "  -0.79%  "

$ws.Range("D18").Value = "'51.523.22"
$ws.Range("E18").Value = "  -0.70%  "

$ws.Range("E19").Value = "  -2.03%  "

$ws.Range("D20").Value = "'7.29"
$ws.Range("E20").Value = "  -4.05%  "

$ws.Range("D21").Value = "'13.15"
$ws.Range("E21").Value = "  -6.41%  "

$ws.Range("D22").Value = "'0.0₃0953"
$ws.Range("E22").Value = "  -2.95%  "

$ws.Range("D23").Value = "'68.90"
$ws.Range("E23").Value = "  -3.00%  "

$ws.Range("D24").Value = "'263.50"
$ws.Range("E24").Value = "  -2.35%  "

$ws.Range("D25").Value = "'2.70"
$ws.Range("E25").Value = "  -4.57%  "

$ws.Range("E26").Value = "  -6.71%  "

$ws.Range("D27").Value = "'26.41"
$ws.Range("E27").Value = "  -2.13%  "

$ws.Range("E28").Value = "  +0.04%  "

$ws.Range("D29").Value = "'7.14"
$ws.Range("E29").Value = "  -6.22%  "

$ws.Range("E30").Value = "  +0.13%  "

$ws.Range("D31").Value = "'6.22"
$ws.Range("E31").Value = "  +2.64%  "

$ws.Range("D32").Value = "'10.03"
$ws.Range("E32").Value = "  -5.24%  "

$ws.Range("E33").Value = "  -0.81%  "

$ws.Range("D34").Value = "'35.33"
$ws.Range("E34").Value = "  -7.20%  "

$ws.Range("D35").Value = "'50.80"
$ws.Range("E35").Value = "  -3.21%  "

$ws.Range("D36").Value = "'1.00"
$ws.Range("E36").Value = "  +0.35%  "

$ws.Range("E37").Value = "  -4.44%  "

$ws.Range("E38").Value = "  -1.01%  "

$ws.Range("D39").Value = "'2.81"
$ws.Range("E39").Value = "  +2.59%  "

$ws.Range("D40").Value = "'17.17"
$ws.Range("E40").Value = "  -6.49%  "

$ws.Range("E41").Value = "  -5.84%  "

$ws.Range("D42").Value = "'0.115"
$ws.Range("E42").Value = "  -4.32%  "

$ws.Range("D43").Value = "'22.77"
$ws.Range("E43").Value = "  -0.92%  "

$ws.Range("D44").Value = "'120.43"
$ws.Range("E44").Value = "  +1.04%  "

$ws.Range("D45").Value = "'2.14"
$ws.Range("E45").Value = "  -1.91%  "

$ws.Range("D46").Value = "'2.087.65"
$ws.Range("E46").Value = "  -1.90%  "

$ws.Range("D47").Value = "'3.23"
$ws.Range("E47").Value = "  -7.54%  "

$ws.Range("E48").Value = "  -6.71%  "

$ws.Range("D49").Value = "'3.228.45"
$ws.Range("E49").Value = "  +0.52%  "

$ws.Range("E50").Value = "  -5.10%  "

$ws.Range("E51").Value = "  -4.91%  "
